$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 419, shifting the existing rows 419-439 down to 420-440
$ws.Rows(419).Insert()

# Populate the newly inserted row 419 with the new weekly data point
$ws.Cells.Item(419, 1).Value = 8
$ws.Cells.Item(419, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(419, 3).Value = "Coquimbo"
$ws.Cells.Item(419, 4).Value = 45075
$ws.Cells.Item(419, 5).Value = 4
$ws.Cells.Item(419, 6).Value = 100112021
$ws.Cells.Item(419, 7).Value = "Ají"
$ws.Cells.Item(419, 8).Value = "Inferno"
$ws.Cells.Item(419, 9).Value = "Primera"
$ws.Cells.Item(419, 10).Value = 360
$ws.Cells.Item(419, 11).Value = 17000
$ws.Cells.Item(419, 12).Value = 18000
$ws.Cells.Item(419, 13).Value = 17500
$ws.Cells.Item(419, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(419, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(419, 16).Value = 1167
$ws.Cells.Item(419, 17).Value = 15
$ws.Cells.Item(419, 18).Value = "Hortaliza"
